# Applies the cryptos list price/volume update described in the commit diff.
# For each changed cell, set the new text value. Columns D that contain
# numeric-looking strings are forced to Text format first so Excel does not
# silently convert them into numbers (which would corrupt formats like
# "0.07137", "1.006", leading/trailing zeros, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.876.20"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.842.62"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.97"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4679"
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3673"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07137"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9334"
$ws.Range("E10").Value = "  +5.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.57"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07689"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "1.858.76"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.282"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.388"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.04"
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008623"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "26.914.52"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.36"
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.017"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.932"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.23"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.22"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.022"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.11"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.883"
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08853"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.202"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.801"
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.170"
$ws.Range("E33").Value = "  +5.66%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7457"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.457"
$ws.Range("E35").Value = "  +0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.082"
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.964"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05184"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5207"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.897"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1515"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.120"
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.52"
$ws.Range("E44").Value = "  +5.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4693"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.006"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "100.52"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.604"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.69"
$ws.Range("E49").Value = "  +1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06027"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8908"
$ws.Range("E51").Value = "  +5.55%  "
